# Junction_Flooding_194: reduce row-5 readings to 2-decimal "custom accuracy"
# and drop the now-redundant row 6 (trimming dimension from A1:AH6 to A1:AH5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (rounded-to-2-decimals) values for B5:AH5, in column order.
$vals = @(
    4.7, 3.37, 0.71, 10.44, 8.01, 3.64, 16.46, 5.82, 2.46, 3.45,
    4.17, 4.52, 1.21, 3.76, 5.25, 3.41, 0.66, 0.37, 49.91, 10.62,
    3.47, 6.91, 3.58, 0.79, 7.84, 3.07, 2.84, 3.33, 4.35, 0.54,
    15.24, 1.84, 4.34
)

for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(5, 2 + $i).Value = $vals[$i]
}

# Remove row 6 entirely (shifts nothing below it up, it was the last row).
$ws.Rows(6).Delete()
